$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume 1h (E) columns for rows 2-48 and 51 ---
# Row 2
$ws.Range("D2").Value = "43.704.44"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value = "2.280.17"
$ws.Range("E3").Value = "  -0.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "124.13"
$ws.Range("E5").Value = "  +8.26%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.12"
$ws.Range("E6").Value = "  -1.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.638"
$ws.Range("E7").Value = "  +1.96%  "

# Row 8
$ws.Range("E8").Value = "  +0.28%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  +1.49%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.29"
$ws.Range("E10").Value = "  -1.54%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0945"
$ws.Range("E11").Value = "  +0.23%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.24"
$ws.Range("E12").Value = "  +2.70%  "

# Row 13
$ws.Range("E13").Value = "  -1.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.43"
$ws.Range("E14").Value = "  -2.51%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.902"
$ws.Range("E15").Value = "  +3.46%  "

# Row 16
$ws.Range("D16").Value = "2.624.82"
$ws.Range("E16").Value = "  -0.43%  "

# Row 17
$ws.Range("D17").Value = "2.277.68"
$ws.Range("E17").Value = "  -0.54%  "

# Row 18
$ws.Range("D18").Value = "43.691.87"
$ws.Range("E18").Value = "  +0.05%  "

# Row 19
$ws.Range("E19").Value = "  +0.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("E20").Value = "  +0.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.38"
$ws.Range("E21").Value = "  +0.13%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.44"
$ws.Range("E22").Value = "  +0.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.66"
$ws.Range("E23").Value = "  +1.06%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.87"
$ws.Range("E24").Value = "  -1.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.43"
$ws.Range("E25").Value = "  -4.07%  "

# Row 26
$ws.Range("E26").Value = "  +1.86%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.82"
$ws.Range("E27").Value = "  +0.81%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.71"
$ws.Range("E28").Value = "  -1.54%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.36"
$ws.Range("E29").Value = "  -0.60%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").Value = "  -0.16%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.67"
$ws.Range("E31").Value = "  -0.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.70"
$ws.Range("E32").Value = "  +0.37%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0915"
$ws.Range("E33").Value = "  -2.25%  "

# Row 34
$ws.Range("E34").Value = "  +1.30%  "

# Row 35
$ws.Range("E35").Value = "  +1.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.26"
$ws.Range("E36").Value = "  +11.56%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0377"
$ws.Range("E37").Value = "  +4.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.65"
$ws.Range("E38").Value = "  -3.24%  "

# Row 39
$ws.Range("E39").Value = "  +0.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("E40").Value = "  +5.79%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.05"
$ws.Range("E41").Value = "  +0.40%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.98"
$ws.Range("E42").Value = "  -5.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.240"
$ws.Range("E43").Value = "  -0.67%  "

# Row 44
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("E45").Value = "  -2.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.61"
$ws.Range("E46").Value = "  -11.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.61"
$ws.Range("E47").Value = "  -2.08%  "

# Row 48
$ws.Range("E48").Value = "  +0.03%  "

# Row 49 (now Cronos; previously ordi)
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").Value = "  +0.31%  "

# Row 50 (now ordi; previously Cronos)
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.60"
$ws.Range("E50").Value = "  +34.33%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.88"
$ws.Range("E51").Value = "  -1.14%  "
